# Apply updated cryptocurrency price/volume data as described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'42.683.92"
$ws.Range('E2').Value = '  -0.01%  '

# Row 3
$ws.Range('D3').Value = "'2.555.76"
$ws.Range('E3').Value = '  +0.85%  '

# Row 4
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.32%  '

# Row 5
$ws.Range('D5').Value = "'311.85"
$ws.Range('E5').Value = '  -1.28%  '

# Row 6
$ws.Range('D6').Value = "'101.11"
$ws.Range('E6').Value = '  +5.03%  '

# Row 7
$ws.Range('D7').Value = "'0.569"
$ws.Range('E7').Value = '  -0.61%  '

# Row 8
$ws.Range('E8').Value = '  +0.30%  '

# Row 9
$ws.Range('D9').Value = "'0.529"
$ws.Range('E9').Value = '  -0.99%  '

# Row 10
$ws.Range('D10').Value = "'36.22"
$ws.Range('E10').Value = '  +2.22%  '

# Row 11
$ws.Range('D11').Value = "'0.0802"
$ws.Range('E11').Value = '  -0.66%  '

# Row 12
$ws.Range('D12').Value = "'7.38"
$ws.Range('E12').Value = '  -1.01%  '

# Row 13
$ws.Range('D13').Value = "'0.108"
$ws.Range('E13').Value = '  +0.38%  '

# Row 14
$ws.Range('D14').Value = "'2.951.01"
$ws.Range('E14').Value = '  +1.03%  '

# Row 15
$ws.Range('D15').Value = "'16.09"
$ws.Range('E15').Value = '  +7.36%  '

# Row 16
$ws.Range('D16').Value = "'2.536.22"
$ws.Range('E16').Value = '  -1.09%  '

# Row 17
$ws.Range('D17').Value = "'0.841"
$ws.Range('E17').Value = '  -0.03%  '

# Row 18
$ws.Range('D18').Value = "'42.707.97"
$ws.Range('E18').Value = '  +0.14%  '

# Row 19
$ws.Range('D19').Value = "'6.81"
$ws.Range('E19').Value = '  +0.02%  '

# Row 20
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = "'0.0₃0954"
$ws.Range('E20').Value = '  -0.14%  '

# Row 21
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').Value = "'12.32"
$ws.Range('E21').Value = '  -1.60%  '

# Row 22
$ws.Range('D22').Value = "'68.98"
$ws.Range('E22').Value = '  -0.40%  '

# Row 23
$ws.Range('D23').Value = "'243.05"
$ws.Range('E23').Value = '  -3.88%  '

# Row 24
$ws.Range('D24').Value = "'2.91"
$ws.Range('E24').Value = '  -1.09%  '

# Row 25
$ws.Range('D25').Value = "'2.06"
$ws.Range('E25').Value = '  +0.67%  '

# Row 26
$ws.Range('D26').Value = "'26.41"
$ws.Range('E26').Value = '  -0.47%  '

# Row 27
$ws.Range('E27').Value = '  -0.06%  '

# Row 28
$ws.Range('D28').Value = "'40.63"
$ws.Range('E28').Value = '  +0.87%  '

# Row 29
$ws.Range('D29').Value = "'2.36"
$ws.Range('E29').Value = '  -1.74%  '

# Row 30
$ws.Range('D30').Value = "'10.12"
$ws.Range('E30').Value = '  -0.83%  '

# Row 31
$ws.Range('D31').Value = "'158.27"
$ws.Range('E31').Value = '  +1.23%  '

# Row 32
$ws.Range('D32').Value = "'5.69"
$ws.Range('E32').Value = '  -2.08%  '

# Row 33
$ws.Range('D33').Value = "'2.76"
$ws.Range('E33').Value = '  +15.93%  '

# Row 34
$ws.Range('D34').Value = "'0.0802"
$ws.Range('E34').Value = '  +0.86%  '

# Row 35
$ws.Range('D35').Value = "'2.06"
$ws.Range('E35').Value = '  -0.87%  '

# Row 36
$ws.Range('E36').Value = '  -2.85%  '

# Row 37
$ws.Range('D37').Value = "'3.20"
$ws.Range('E37').Value = '  -2.51%  '

# Row 38
$ws.Range('D38').Value = "'18.23"
$ws.Range('E38').Value = '  -5.98%  '

# Row 39
$ws.Range('E39').Value = '  -1.09%  '

# Row 40
$ws.Range('E40').Value = '  -0.10%  '

# Row 41
$ws.Range('D41').Value = "'4.24"
$ws.Range('E41').Value = '  +11.92%  '

# Row 42
$ws.Range('D42').Value = "'21.55"
$ws.Range('E42').Value = '  -4.52%  '

# Row 43
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = "'1.01"
$ws.Range('E43').Value = '  +0.46%  '

# Row 44
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D44').Value = "'3.33"
$ws.Range('E44').Value = '  +2.99%  '

# Row 45
$ws.Range('E45').Value = '  -1.00%  '

# Row 46
$ws.Range('D46').Value = "'1.967.03"
$ws.Range('E46').Value = '  -1.08%  '

# Row 47
$ws.Range('D47').Value = "'8.97"
$ws.Range('E47').Value = '  -0.18%  '

# Row 48
$ws.Range('D48').Value = "'2.809.79"
$ws.Range('E48').Value = '  +1.73%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'0.193"
$ws.Range('E49').Value = '  +1.39%  '

# Row 50
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').Value = "'80.93"
$ws.Range('E50').Value = '  -3.80%  '

# Row 51
$ws.Range('D51').Value = "'73.16"
$ws.Range('E51').Value = '  -1.23%  '
